$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 3047.5833
$ws.Range("I96").Value = 685.6667
$ws.Range("J96").Value = 10133.333
$ws.Range("K96").Value = 2057.0001
$ws.Range("L96").Value = 30399.999
$ws.Range("M96").Value = -684.0001000000002
$ws.Range("N96").Value = -33145.999

$ws.Range("H98").Value = 6234405.5
$ws.Range("I98").Value = 73745.6
$ws.Range("K98").Value = 73745.6
$ws.Range("M98").Value = -72247.6

$ws.Range("H107").Value = 8258.267
$ws.Range("I107").Value = 9544.583
$ws.Range("J107").Value = 3113
$ws.Range("K107").Value = 9544.583
$ws.Range("L107").Value = 3113
$ws.Range("M107").Value = -7624.583000000001
$ws.Range("N107").Value = -6953

$ws.Range("H122").Value = 6234405.5
$ws.Range("I122").Value = 73745.6
$ws.Range("K122").Value = 221236.8
$ws.Range("M122").Value = -218786.8

$ws.Range("H140").Value = 85997.5
$ws.Range("J140").Value = 85997.5
$ws.Range("L140").Value = 85997.5
$ws.Range("N140").Value = -96357.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2356.4285
$ws.Range("I97").Value = 1980.4375
$ws.Range("J97").Value = 3559.6
$ws.Range("K97").Value = 1980.4375
$ws.Range("L97").Value = 3559.6
$ws.Range("M97").Value = -1484.4375
$ws.Range("N97").Value = -4551.6

$ws.Range("H122").Value = 3916.7932
$ws.Range("I122").Value = 4151.7393
$ws.Range("J122").Value = 3016.1667
$ws.Range("K122").Value = 12455.2179
$ws.Range("L122").Value = 9048.500100000001
$ws.Range("M122").Value = -10005.2179
$ws.Range("N122").Value = -13948.5001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 871.5294
$ws.Range("I94").Value = 602.53845
$ws.Range("J94").Value = 1745.75
$ws.Range("K94").Value = 602.53845
$ws.Range("L94").Value = 1745.75
$ws.Range("M94").Value = -151.53845
$ws.Range("N94").Value = -2647.75

$ws.Range("H107").Value = 228089.61
$ws.Range("I107").Value = 747.625
$ws.Range("J107").Value = 834334.94
$ws.Range("K107").Value = 747.625
$ws.Range("L107").Value = 834334.94
$ws.Range("M107").Value = 1172.375
$ws.Range("N107").Value = -838174.94

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1466.6666
$ws.Range("I117").Value = 840
$ws.Range("J117").Value = 2250
$ws.Range("K117").Value = 2520
$ws.Range("L117").Value = 6750
$ws.Range("M117").Value = 922
$ws.Range("N117").Value = -13634

$ws.Range("H126").Value = 2340
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 2550
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 7650
$ws.Range("M126").Value = 440
$ws.Range("N126").Value = -17530

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 382809.53
$ws.Range("J33").Value = 382809.53
$ws.Range("L33").Value = 382809.53
$ws.Range("N33").Value = -383313.53

$ws.Range("H97").Value = 1772.45
$ws.Range("I97").Value = 1107.8572
$ws.Range("J97").Value = 3323.1667
$ws.Range("K97").Value = 1107.8572
$ws.Range("L97").Value = 3323.1667
$ws.Range("M97").Value = -611.8571999999999
$ws.Range("N97").Value = -4315.1667

$ws.Range("H102").Value = 1734.3462
$ws.Range("I102").Value = 1804.5264
$ws.Range("J102").Value = 1543.8572
$ws.Range("K102").Value = 1804.5264
$ws.Range("L102").Value = 1543.8572
$ws.Range("M102").Value = -182.5264
$ws.Range("N102").Value = -4787.8572

$ws.Range("H122").Value = 2658.7646
$ws.Range("I122").Value = 2468
$ws.Range("J122").Value = 3278.75
$ws.Range("K122").Value = 7404
$ws.Range("L122").Value = 9836.25
$ws.Range("M122").Value = -4954
$ws.Range("N122").Value = -14736.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 548.5
$ws.Range("I22").Value = 534.5769
$ws.Range("J22").Value = 593.75
$ws.Range("K22").Value = 534.5769
$ws.Range("L22").Value = 593.75
$ws.Range("M22").Value = -239.5769
$ws.Range("N22").Value = -1183.75

$ws.Range("H27").Value = 548.5
$ws.Range("I27").Value = 534.5769
$ws.Range("J27").Value = 593.75
$ws.Range("K27").Value = 534.5769
$ws.Range("L27").Value = 593.75
$ws.Range("M27").Value = -427.5769
$ws.Range("N27").Value = -807.75

$ws.Range("H93").Value = 1698.3182
$ws.Range("I93").Value = 1164.8572
$ws.Range("J93").Value = 2631.875
$ws.Range("K93").Value = 1164.8572
$ws.Range("L93").Value = 2631.875
$ws.Range("M93").Value = 83.14280000000008
$ws.Range("N93").Value = -5127.875

$ws.Range("H100").Value = 2674.9546
$ws.Range("I100").Value = 1845.5834
$ws.Range("J100").Value = 3670.2
$ws.Range("K100").Value = 1845.5834
$ws.Range("L100").Value = 3670.2
$ws.Range("M100").Value = -1304.5834
$ws.Range("N100").Value = -4752.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H54").Value = 6500
$ws.Range("J54").Value = 6500
$ws.Range("L54").Value = 6500
$ws.Range("N54").Value = -7540

$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H61").Value = 8500
$ws.Range("I61").Value = 7000
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 7000
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -6708
$ws.Range("N61").Value = -10584

$ws.Range("H96").Value = 3476.923
$ws.Range("I96").Value = 2447.0588
$ws.Range("J96").Value = 5422.222
$ws.Range("K96").Value = 2447.0588
$ws.Range("L96").Value = 5422.222
$ws.Range("M96").Value = -1074.0588
$ws.Range("N96").Value = -8168.222

$ws.Range("H100").Value = 1064.5
$ws.Range("I100").Value = 1016.1111
$ws.Range("K100").Value = 2032.2222
$ws.Range("M100").Value = -1491.2222
